# B1-and-B2-PowerPoint.pptx edit
#
# 1. Slide 5's table switches from the deck's custom "Table_0" style
#    ({A798E453-6D37-4A86-98C8-6696E5F8EF2D}, defined in ppt/tableStyles.xml)
#    to the built-in "Medium Style 2 - Accent 1" style
#    ({3D92BB1B-A5A2-4353-8444-DBD72E06E5AD}).
#
# 2. The slide-master theme ("Integral" / "Red Violet" colour scheme) and the
#    notes-master theme ("Office Theme" / "Office" colour scheme) swap their
#    colour schemes with one another (font scheme / format scheme are already
#    identical between the two themes, so only the 12 theme colours move).

$p = $ppt.ActivePresentation

function HexToRgbInt($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# ---------------------------------------------------------------------------
# 1. Table style on slide 5
# ---------------------------------------------------------------------------
$slide5 = $p.Slides.Item(5)
$table = $slide5.Shapes.Item(2).Table
$table.ApplyStyle("{3D92BB1B-A5A2-4353-8444-DBD72E06E5AD}")

# ---------------------------------------------------------------------------
# 2. Swap the two theme colour schemes
# ---------------------------------------------------------------------------
# Order used by ThemeColorScheme.Colors(): dk1, lt1, dk2, lt2, accent1..6,
# hlink, folHlink.
$officeColors = @("000000", "FFFFFF", "44546A", "E7E6E6", "5B9BD5", "ED7D31", "A5A5A5", "FFC000", "4472C4", "70AD47", "0563C1", "954F72")
$integralColors = @("000000", "FFFFFF", "454551", "D8D9DC", "E32D91", "C830CC", "4EA6DC", "4775E7", "8971E1", "D54773", "6B9F25", "8C8C8C")

$masterScheme = $p.SlideMaster.Theme.ThemeColorScheme
$notesScheme = $p.NotesMaster.Theme.ThemeColorScheme

for ($i = 1; $i -le 12; $i++) {
    # Slide master ("Integral") becomes the Office colours.
    $masterScheme.Colors($i).RGB = HexToRgbInt $officeColors[$i - 1]
    # Notes master ("Office Theme") becomes the Integral/Red Violet colours.
    $notesScheme.Colors($i).RGB = HexToRgbInt $integralColors[$i - 1]
}
